$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.236.76'
$ws.Range("E2").Value = '  +1.36%  '

# Row 3
$ws.Range("D3").Value = '1.801.25'
$ws.Range("E3").Value = '  +2.69%  '

# Row 4
$ws.Range("E4").Value = '  -0.45%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '339.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.29%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.11%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4671'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +24.41%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3637'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.82%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.48'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.87%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.144'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.59%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07586'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.69%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.51'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.01%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.001'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.40%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.254'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.76%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.276'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.22%  '

# Row 16
$ws.Range("D16").Value = '1.796.00'
$ws.Range("E16").Value = '  +2.06%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001088'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.16%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06703'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.32%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.10%  '

# Row 20
$ws.Range("E20").Value = '  -0.31%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.76%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.398'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.60%  '

# Row 23
$ws.Range("D23").Value = '28.230.28'
$ws.Range("E23").Value = '  +1.18%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.89%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.406'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.40%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.48'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.52%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.402'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.43%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.24'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.20%  '

# Row 29
$ws.Range("D29").Value = '1.999.31'
$ws.Range("E29").Value = '  +1.97%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.275'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.67%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '133.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.40%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.070'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.34%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.905'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.11%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09535'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.43%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02373'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.28%  '

# Row 36
$ws.Range("E36").Value = '  +0.20%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06287'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.20%  '

# Row 38
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6642'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.94%  '

# Row 39
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.201'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.06%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2169'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.54%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.479'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.35%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.215'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.75%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.084'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.30%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9987'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.30%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.60%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.871'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.18%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6096'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.96%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.54'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.14%  '

# Row 49
$ws.Range("E49").Value = '  +2.04%  '

# Row 50
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.170'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.14%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07092'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.43%  '
